$d = $word.ActiveDocument

# --- Change 1: split the sentence run and insert a _GoBack bookmark ---
# Find the location right before " más vendidos." so the bookmark sits between
# "...productos" and " más vendidos."
$splitPoint = $d.Content
$splitPoint.Find.Execute(" más vendidos.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($splitPoint.Start, $splitPoint.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Change 2: remove the "Especificación de Requerimientos..." section,
# the "Historias de Usuario" table, and the "Tabla 1..." caption, keeping a
# single blank paragraph in their place ---

# Locate the blank paragraph that should be kept (the first empty paragraph
# right after the "más vendidos." sentence) and find its index. The
# paragraph immediately after the sentence starts right where the sentence's
# own paragraph mark ends (Find's End does not include that mark, hence +1).
$afterSentence = $d.Content
$afterSentence.Find.Execute("más vendidos.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$keepIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq ($afterSentence.End + 1)) {
        $keepIndex = $i
        break
    }
}

# Delete the 19 paragraphs (blank + headings) that precede the table.
$deleteIndex = $keepIndex + 1
for ($k = 1; $k -le 19; $k++) {
    $d.Paragraphs.Item($deleteIndex).Range.Delete()
}

# Delete the 4 paragraphs that remain before the table itself
# ("Especificación de Requerimientos del software", blank,
# "Historias de Usuario:", blank).
for ($k = 1; $k -le 4; $k++) {
    $d.Paragraphs.Item($deleteIndex).Range.Delete()
}

# Delete the "Historias de Usuario" table entirely.
$d.Tables.Item(1).Delete()

# Delete the "Tabla 1 - Requisito funcional 1" caption paragraph (text +
# its trailing paragraph mark) using Find, since the Paragraphs collection
# is not reliable immediately after a table deletion.
$caption = $d.Content
$caption.Find.Execute("Tabla 1 - Requisito funcional 1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Range($caption.Start, $caption.End + 1).Delete()
